# Add 8 new "Title and Content" slides (ppLayoutText = 2) carrying the
# GCD / XCon paper summary, matching the target deck.
$p = $ppt.ActivePresentation

$slidesData = @(
    @{
        Title = "Introduction to Generalized Category Discovery (GCD) and XCon"
        Bullets = @(
            "- GCD: a practical approach for categorizing unlabeled instances in the absence of knowledge about all classes",
            "- XCon: a method for fine-grained GCD using expert datasets and contrastive learning"
        )
    },
    @{
        Title = "Motivation: Challenges in Fine-Grained Category Discovery"
        Bullets = @(
            "- Large inter-class similarities and intra-class variances in fine-grained datasets",
            "- Traditional methods often struggle with irrelevant cues, leading to suboptimal performance"
        )
    },
    @{
        Title = "Related Work: Novel Category Discovery (NCD) and Contrastive Learning"
        Bullets = @(
            "- NCD: categorizing unseen classes based on knowledge from seen classes",
            "- Transfer learning and self-supervision techniques used in NCD",
            "- Limitations of contrastive learning for fine-grained classification"
        )
    },
    @{
        Title = "Methodology: XCon's Approach to Fine-Grained Category Discovery"
        Bullets = @(
            "- Partitioning dataset into expert datasets using k-means clustering on self-supervised representations",
            "- Learning from both coarse-grained and fine-grained features through supervised and unsupervised contrastive learning"
        )
    },
    @{
        Title = "Experiments and Results: State-of-the-Art Performance on Various Datasets"
        Bullets = @(
            "- CIFAR-10/100, ImageNet-100, CUB-200, Stanford Cars, FGVC-Aircraft, Oxford-IIIT Pet",
            "- Achieved state-of-the-art performance on several fine-grained category discovery benchmarks"
        )
    },
    @{
        Title = "Implementation Details: Model Architecture and Training"
        Bullets = @(
            "- ViT-B-16 model initialized with DINO-pretrained parameters",
            "- Fine-tuning final transformer block",
            "- Combination of supervised and unsupervised contrastive losses"
        )
    },
    @{
        Title = "Ablation Studies: Validating XCon's Effectiveness"
        Bullets = @(
            "- Impact of fine-grained loss weight and number of sub-datasets on performance",
            "- Robustness of XCon across different configurations"
        )
    },
    @{
        Title = "Conclusion: Significant Advancement in Fine-Grained Category Discovery"
        Bullets = @(
            "- XCon's success on various benchmarks marks a significant advancement in the field",
            "- Encourages further exploration and adaptation within the research community."
        )
    }
)

$idx = 1
foreach ($slideData in $slidesData) {
    # ppLayoutText (2) => "Title and Content" layout (title placeholder + body placeholder)
    $s = $p.Slides.Add($idx, 2)

    $titleShape = $s.Shapes.Item(1)
    [void]$titleShape.TextFrame.TextRange.InsertAfter($slideData.Title)

    $bodyShape = $s.Shapes.Item(2)
    $bodyText = [string]::Join("`r", $slideData.Bullets)
    [void]$bodyShape.TextFrame.TextRange.InsertAfter($bodyText)

    $idx = $idx + 1
}
